{"js": "// Load all paragraphs in the document body.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// 1) Insert a new empty paragraph right after the paragraph that contains\n//    \"You need to build an AI Agent for\" (the one with the hyperlink) and\n//    right before the \"Here is the gdrive link of pre-scraped data.\" paragraph.\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"You need to build an AI Agent for\") !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\nif (anchor) {\n  anchor.insertParagraph(\"\", Word.InsertLocation.after);\n  await context.sync();\n}\n\n// 2) & 3) Update the \"Scoring rules\" bullet text.\nconst body = context.document.body;\n\nconst find1 = body.search(\n  \"Proper Human-in-the-loop Usage (Agent should be able to get feedback from humans based on the situation. For example if the question of the human is not clear, the agent can ask for human input clarification.) - (point 2)\",\n  { matchCase: true }\n);\nfind1.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < find1.items.length; i++) {\n  find1.items[i].insertText(\n    \"Proper Human-in-the-loop Usage (E.g. Agents should be able to get feedback from humans based on the situation. For example if the question of the human is not clear, the agent can ask for human input clarification.) - (point 2)\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\nconst find2 = body.search(\n  \"Reasoning Capability (This should be implemented as an AI reasoning agent. For development, please use non-reasoning models.) - (point 3)\",\n  { matchCase: true }\n);\nfind2.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < find2.items.length; i++) {\n  find2.items[i].insertText(\n    \"Reasoning Capability (E.g. This should be implemented as an AI reasoning agent. For development, please use non-reasoning models. And reasoning steps and intermediate results should be displayed in the frontend) - (point 3)\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Insert a new empty paragraph right after the paragraph that contains\n#    \"You need to build an AI Agent for\" (the hyperlink paragraph) and right\n#    before the \"Here is the gdrive link of pre-scraped data.\" paragraph.\nforeach ($p in $d.Paragraphs) {\n  if ($p.Range.Text -like \"*You need to build an AI Agent for*\") {\n    $p.Range.InsertParagraphAfter()\n    break\n  }\n}\n\n# 2) Update the \"Human-in-the-loop\" bullet text.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Proper Human-in-the-loop Usage (Agent should be able to get feedback from humans based on the situation. For example if the question of the human is not clear, the agent can ask for human input clarification.) - (point 2)\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Proper Human-in-the-loop Usage (E.g. Agents should be able to get feedback from humans based on the situation. For example if the question of the human is not clear, the agent can ask for human input clarification.) - (point 2)\"\n$find.Execute([ref]$find.Text, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]0, [ref]$false, [ref]$find.Replacement.Text, [ref]2)\n\n# 3) Update the \"Reasoning Capability\" bullet text.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"Reasoning Capability (This should be implemented as an AI reasoning agent. For development, please use non-reasoning models.) - (point 3)\"\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"Reasoning Capability (E.g. This should be implemented as an AI reasoning agent. For development, please use non-reasoning models. And reasoning steps and intermediate results should be displayed in the frontend) - (point 3)\"\n$find2.Execute([ref]$find2.Text, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]0, [ref]$false, [ref]$find2.Replacement.Text, [ref]2)\n\n$d.Save()\n"}
